$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing row 253 (low/close were revised)
$ws.Range("E253").Value = 1.77126
$ws.Range("F253").Value = 1.7811

# New monthly rows appended after row 253, matching the source format
$newRows = @(
    @{ Row = 254; A = 45047.33333333334; C = 1.7837;  D = 1.83079; E = 1.76615; F = 1.83079; G = 0 },
    @{ Row = 255; A = 45078.33333333334; C = 1.82839; D = 1.83079; E = 1.78046; F = 1.79995; G = 0 },
    @{ Row = 256; A = 45110.33333333334; C = 1.7945;  D = 1.7978;  E = 1.7945;  F = 1.79632; G = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the formatting (style index) of the last data row's date cell
    # down into the new date cell before setting its value.
    $ws.Range("A253").Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "FX_IDC:USDBGN"
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
